$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1) ---
$ws = $wb.Worksheets.Item(1)

# Version bump 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date bump
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher gains a value
$ws.Range("B9").Value = "Alvearie Team"

# The duplicated "Contact" / "No display for ContactDetail" row (row 10) becomes "Jurisdiction"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the now-redundant duplicate "Contact" row (old row 11)
$ws.Rows.Item(11).Delete()

# --- Sheet "Elements" (sheet2) ---
$ws2 = $wb.Worksheets.Item(2)

# Root "Extension" element's Short/Definition describe the renamed resource
$ws2.Range("K2").Value = "RelatedIssue"
$ws2.Range("L2").Value = "Related issues that can be combined to fulfill a single time period for a given care gap."
